$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric values (e.g. "1.035") are kept as text, not converted to numbers
$ws.Range("B2:E51").NumberFormat = "@"

$data = @(
    2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '28.830.66', '  +3.09%  ',
    3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.915.53', '  +3.17%  ',
    4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.035', '  +3.07%  ',
    5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '322.06', '  +3.13%  ',
    6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.034', '  +3.11%  ',
    7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5230', '  +1.84%  ',
    8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3964', '  +3.60%  ',
    9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08395', '  +2.17%  ',
    10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.141', '  +2.94%  ',
    11, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '42.73', '  +3.08%  ',
    12, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.333', '  +2.59%  ',
    13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.920.66', '  +2.88%  ',
    14, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.77', '  +1.42%  ',
    15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.350', '  +1.41%  ',
    16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.036', '  +3.11%  ',
    17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001117', '  +1.98%  ',
    18, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '92.12', '  +1.89%  ',
    19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06846', '  +3.00%  ',
    20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '18.05', '  +2.32%  ',
    21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.033', '  +2.95%  ',
    22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.122', '  +1.98%  ',
    23, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '28.863.26', '  +3.06%  ',
    24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.32', '  +2.64%  ',
    25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.275', '  +0.78%  ',
    26, 'LEO', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', '3.430', '  +1.21%  ',
    27, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.125.29', '  +2.52%  ',
    28, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '163.50', '  +4.12%  ',
    29, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '21.09', '  +3.37%  ',
    30, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.462', '  -1.58%  ',
    31, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '128.20', '  +2.97%  ',
    32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1063', '  -0.13%  ',
    33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.061', '  +3.12%  ',
    34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.022', '  +2.06%  ',
    35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.702', '  +2.90%  ',
    36, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.531', '  +1.66%  ',
    37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06700', '  +2.95%  ',
    38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02483', '  +3.05%  ',
    39, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2236', '  +2.78%  ',
    40, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6605', '  +1.18%  ',
    41, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.273', '  +5.20%  ',
    42, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.200', '  +0.89%  ',
    43, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.059', '  +1.78%  ',
    44, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.23', '  +0.84%  ',
    45, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6198', '  +0.95%  ',
    46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '13.25', '  +1.80%  ',
    47, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.779', '  +2.99%  ',
    48, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.309', '  +2.42%  ',
    49, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.248', '  +3.08%  ',
    50, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.028', '  +1.03%  ',
    51, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '123.47', '  +2.20%  '
)

for ($i = 0; $i -lt $data.Count; $i = $i + 5) {
    $r = $data[$i]
    $ws.Cells.Item($r, 2).Value = $data[$i + 1]
    $ws.Cells.Item($r, 3).Value = $data[$i + 2]
    $ws.Cells.Item($r, 4).Value = $data[$i + 3]
    $ws.Cells.Item($r, 5).Value = $data[$i + 4]
}
